$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (A:C) ---------------------------------------------------
# Target OOXML width is 37.140625; closest value reachable through the
# ColumnWidth property of this runtime is ~37.1667 (internal width is
# quantized), so use an input value that lands on that nearest value.
$ws.Columns.Item(1).ColumnWidth = 37.140625
$ws.Columns.Item(2).ColumnWidth = 37.140625
$ws.Columns.Item(3).ColumnWidth = 37.140625

# --- Row heights -------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 55.5
$ws.Rows.Item(2).RowHeight = 13.5
$ws.Rows.Item(3).RowHeight = 13.5

# --- New column N: copy formatting from column M for rows 3, 4 and 8 ---
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)

$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 2023

$ws.Range("N5").Value = 4.3499999999999996
$ws.Range("N6").Value = 4.3499999999999996

$ws.Range("N7").Value = "-"
$ws.Range("N7").HorizontalAlignment = -4152
$ws.Range("N7").VerticalAlignment = -4108
$ws.Range("N7").Font.Name = "Times New Roman"
$ws.Range("N7").Font.Size = 9

$ws.Range("M8").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N8").Value = "-"

# --- Reset the active cell / selection to the default top-left cell ----
$ws.Range("A1").Select()
